# Set assignees_id default value to the migration user for the "warnings"
# entity across every mapping sheet that defines it.
#
# Each of these sheets has, on row 3 (the "added_by" / assignees:id mapping
# row):
#   column M = default_value
#   column P = is_complete
#   column Q = comments
# The migration user id (2657) is now used as the default value, the row is
# marked complete ("yes") and the comment explains the decision.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "client_violent_warnings",
    "p1_client_remarks_warnings",
    "deputy_violent_warnings",
    "client_special_warnings",
    "deputy_special_warnings",
    "client_saarcheck_warnings",
    "client_nodebtchase_warnings"
)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("M3").Value = 2657
    $ws.Range("P3").Value = "yes"
    $ws.Range("Q3").Value = "Set to migration user"
}
